$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2
Set-TextValue "D2" "293.27"
Set-TextValue "E2" "-4.87%"

# Row 3
Set-TextValue "D3" "40.58"
Set-TextValue "E3" "-1.17%"

# Row 4
Set-TextValue "D4" "5.028"
Set-TextValue "E4" "-2.31%"

# Row 5
Set-TextValue "D5" "0.07330"
Set-TextValue "E5" "-3.64%"

# Row 6
Set-TextValue "D6" "1.542"
Set-TextValue "E6" "-8.10%"

# Row 7
Set-TextValue "E7" "-0.72%"

# Row 8
Set-TextValue "D8" "2.360"
Set-TextValue "E8" "-2.64%"

# Row 9
Set-TextValue "D9" "0.1165"
Set-TextValue "E9" "-2.86%"

# Row 10
Set-TextValue "D10" "0.1751"
Set-TextValue "E10" "-3.83%"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.08749"
Set-TextValue "E11" "-3.31%"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.04340"
Set-TextValue "E12" "5.06%"

# Row 13
Set-TextValue "D13" "0.1056"
Set-TextValue "E13" "0.17%"

# Row 14
Set-TextValue "D14" "0.001276"
Set-TextValue "E14" "-1.50%"

# Row 15
Set-TextValue "D15" "0.005951"
Set-TextValue "E15" "2.05%"

# Row 16
Set-TextValue "D16" "3.340"
Set-TextValue "E16" "0.02%"

# Row 17
Set-TextValue "E17" "-0.80%"

# Row 18
Set-TextValue "E18" "-1.98%"

# Row 19
Set-TextValue "D19" "7.974"
Set-TextValue "E19" "5.04%"

# Row 20
Set-TextValue "D20" "0.1390"
Set-TextValue "E20" "3.63%"

# Row 22
Set-TextValue "D22" "0.03933"
Set-TextValue "E22" "-0.53%"

# Row 23
Set-TextValue "D23" "0.001262"
Set-TextValue "E23" "-1.40%"

# Row 24
Set-TextValue "D24" "0.003678"
Set-TextValue "E24" "-9.22%"

# Row 25
Set-TextValue "E25" "-5.21%"

# Row 26
Set-TextValue "E26" "22.44%"

# Row 38
Set-TextValue "D38" "0.02318"
Set-TextValue "E38" "-4.31%"

# Row 39
Set-TextValue "D39" "0.05070"
Set-TextValue "E39" "-1.87%"

# Row 40
Set-TextValue "D40" "0.006196"
Set-TextValue "E40" "87.61%"

# Row 41
Set-TextValue "D41" "0.007847"
Set-TextValue "E41" "1.64%"

# Row 42
Set-TextValue "D42" "0.1287"
Set-TextValue "E42" "-1.14%"

# Row 43
Set-TextValue "D43" "0.007372"
Set-TextValue "E43" "-2.97%"

# Row 44
Set-TextValue "D44" "0.007260"
Set-TextValue "E44" "-4.17%"

# Row 45
Set-TextValue "D45" "0.3193"
Set-TextValue "E45" "-4.06%"

# Row 46
Set-TextValue "D46" "0.00006271"
Set-TextValue "E46" "-5.16%"

# Row 47
Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "0.02%"

# Row 48
Set-TextValue "E48" "-87.71%"

# Row 49
Set-TextValue "D49" "0.00002102"
Set-TextValue "E49" "0.02%"

# Row 50
Set-TextValue "D50" "0.0002002"
Set-TextValue "E50" "0.02%"
